# Insert a new product row at row 6 of the "Artículos" sheet, pushing the
# existing rows (old row 6 onward) down by one, and populate the new row
# with the "click & roll" article data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# Insert a new blank row above the current row 6 (shifts rows 6.. down by 1).
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the article's data.
$ws.Cells.Item(6, 1).Value = 78421974
$ws.Cells.Item(6, 2).Value = "Cigarrillos"
$ws.Cells.Item(6, 3).Value = "rubios"
$ws.Cells.Item(6, 4).Value = "click & roll"
$ws.Cells.Item(6, 5).Value = "Lucky Strike"
$ws.Cells.Item(6, 6).Value = 20
$ws.Cells.Item(6, 7).Value = "und."
$ws.Cells.Item(6, 8).Value = "caja"
$ws.Cells.Item(6, 9).Value = "Cigarrillos"
$ws.Cells.Item(6, 10).Value = "Paraguay"
$ws.Cells.Item(6, 11).Value = 10
$ws.Cells.Item(6, 12).Value = $false
$ws.Cells.Item(6, 13).Value = $true
# Column N (14) intentionally left blank for this row.
$ws.Cells.Item(6, 15).Value = $true
$ws.Cells.Item(6, 16).Value = $true
